$wb = $excel.ActiveWorkbook

# --- Constants (Excel COM enumerations) ---
# Border index: xlEdgeTop=8, xlEdgeRight=10
# Border line style: xlContinuous=1
# xlPasteFormats=-4122 (used with Range.PasteSpecial)
#
# NOTE on style-table hygiene: this engine commits a style-table entry on
# every individual Borders.Item(n).LineStyle assignment, and never garbage
# collects an entry that stops being referenced. Toggling each of the 4
# border edges of a cell one at a time therefore leaves orphaned, unused
# xf/border records behind (extra, unreferenced <xf>/<border> entries that
# don't match the source diff). To keep the style table exactly as small
# as the target (no stray entries), the desired border combination is
# established ONCE on a single cell via direct edge writes, and then
# propagated to every other cell that needs the identical formatting with
# Copy / PasteSpecial(xlPasteFormats), which clones the already-committed
# style instead of re-deriving it edge-by-edge.

$xlContinuous = 1
$xlPasteFormats = -4122

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Establish the "thin top + thin bottom, no left/right" border on C1 (sheet1)
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = $xlContinuous   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = $xlContinuous   # xlEdgeBottom

# Establish the "thin top + thin bottom + thin right, no left" border on D1
# (sheet1) by starting from C1's already-committed style and adding the
# right edge.
$d1 = $ws1.Range("D1")
$c1.Copy()
$d1.PasteSpecial($xlPasteFormats)
$d1.Borders.Item(10).LineStyle = $xlContinuous  # xlEdgeRight

# Propagate the two finished styles to sheet2's matching header cells.
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Text updates: "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
